$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold text-formatted numbers (e.g. thousands-dot
# notation like 58.770.24). Force text storage so Excel's automatic
# number coercion does not turn them into numeric values, then restore
# the default cell style so no stray number-format style is left behind.
$dCells = @("D2", "D3", "D5", "D6", "D9", "D11", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D49")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "58.770.24"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "2.585.51"
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "519.86"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "139.28"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "2.594.58"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").Value = "0.100"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D14").Value = "3.038.71"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "58.747.11"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "20.40"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "2.551.73"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").Value = "338.43"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "10.14"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "6.51"
$ws.Range("E22").Value = "  +4.74%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "66.26"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "7.03"
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "0.0₃0718"
$ws.Range("E30").Value = "  -4.04%  "
$ws.Range("D31").Value = "5.93"
$ws.Range("E31").Value = "  -4.96%  "
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "18.77"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "148.83"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "3.96"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "1.46"
$ws.Range("E38").Value = "  +2.10%  "
$ws.Range("D39").Value = "0.827"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "274.80"
$ws.Range("E43").Value = "  +2.29%  "
$ws.Range("D44").Value = "10.74"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").Value = "0.589"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "0.0948"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("D48").Value = "18.51"
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("D49").Value = "1.980.53"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E51").Value = "  -0.92%  "

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
